$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 2000年-2009年 rows (old rows 2-11); this shifts the 2010年-2020年
# data (old rows 12-22) up to rows 2-12 and updates the sheet dimension.
$ws.Rows.Item(2).Resize(10).Delete() | Out-Null

# The 2020年 exchange-rate-to-euro figure was corrected/rounded.
$ws.Range("C12").Value = 787.55

# Append the newly published 2021年 and 2022年 rows.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 5.8735
$ws.Range("C13").Value = 762.9299999999999
$ws.Range("D13").Value = 83
$ws.Range("E13").Value = 645.15

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 5.1261
$ws.Range("C14").Value = 707.21
$ws.Range("D14").Value = 85.89
$ws.Range("E14").Value = 672.61

# Match the existing year-label formatting (bold, centered, bordered).
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13:A14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
